# Generate Report for Handback
# Update the "Latest Handback DateTime" (column K, row 2 - the
# 29bc900b-f4ef-4d73-b187-8d5b003fec25 entry) on both the zh-cn and
# de-de localization-status sheets to reflect the newly generated
# handback report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("K2").Value = "2016-10-19 11:47:49"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-10-19 11:48:08"
